# Replace NA with "" for the NOTEX_da / NOTEX_en / NOTEX_kl rows on the
# General_MD sheet (rows 13-15, column B). The cells previously held the
# placeholder text "Befolkningsstatistikregistret indeholder …"; clearing
# their contents (while keeping the existing cell style/format) removes
# that now-unused shared string from the workbook.

$wb = $excel.ActiveWorkbook

$wsGeneral = $wb.Worksheets.Item("General_MD")

$wsGeneral.Range("B13").ClearContents() | Out-Null
$wsGeneral.Range("B14").ClearContents() | Out-Null
$wsGeneral.Range("B15").ClearContents() | Out-Null

# Make General_MD the active/selected sheet, with B13:B15 selected, matching
# the editor's last interaction before saving.
$wsGeneral.Activate() | Out-Null
$wsGeneral.Range("B13:B15").Select() | Out-Null
